$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 1.233604490811956
$ws.Cells.Item(2, 3).Value = 0.07243487903528489
$ws.Cells.Item(2, 4).Value = 0.004012937065372757
$ws.Cells.Item(2, 5).Value = 0.0651270280068692
$ws.Cells.Item(2, 6).Value = 4.639663493202576
$ws.Cells.Item(2, 8).Value = 0.07973214163530429
$ws.Cells.Item(2, 9).Value = 3.007647047798784
$ws.Cells.Item(2, 10).Value = 0.1738394188915784
$ws.Cells.Item(2, 11).Value = 1.098044654369545
$ws.Cells.Item(2, 12).Value = 0.3540689197455009

$ws.Cells.Item(3, 2).Value = 1.216754902809356
$ws.Cells.Item(3, 3).Value = 0.06715396461611078
$ws.Cells.Item(3, 4).Value = 0.004145253105960034
$ws.Cells.Item(3, 5).Value = 0.06537049419064811
$ws.Cells.Item(3, 6).Value = 4.599691938735035
$ws.Cells.Item(3, 8).Value = 0.07973214163530429
$ws.Cells.Item(3, 9).Value = 2.983223206827091
$ws.Cells.Item(3, 10).Value = 0.173739443769918
$ws.Cells.Item(3, 11).Value = 1.075287720699464
$ws.Cells.Item(3, 12).Value = 0.3524579666636996

$ws.Cells.Item(4, 2).Value = 1.207252263149655
$ws.Cells.Item(4, 3).Value = 0.06395599346494407
$ws.Cells.Item(4, 4).Value = 0.00423337219957487
$ws.Cells.Item(4, 5).Value = 0.06554090196965134
$ws.Cells.Item(4, 6).Value = 4.576408819635773
$ws.Cells.Item(4, 8).Value = 0.07973214163530429
$ws.Cells.Item(4, 9).Value = 2.968916534873586
$ws.Cells.Item(4, 10).Value = 0.1737113378702055
$ws.Cells.Item(4, 11).Value = 1.062085975564344
$ws.Cells.Item(4, 12).Value = 0.3516446564154734

$ws.Cells.Item(5, 2).Value = 1.203592118987615
$ws.Cells.Item(5, 3).Value = 0.06266391640198776
$ws.Cells.Item(5, 4).Value = 0.004271017069480187
$ws.Cells.Item(5, 5).Value = 0.06561561930578197
$ws.Cells.Item(5, 6).Value = 4.567237214676211
$ws.Cells.Item(5, 8).Value = 0.07973214163530429
$ws.Cells.Item(5, 9).Value = 2.96325963005242
$ws.Cells.Item(5, 10).Value = 0.1737082610827834
$ws.Cells.Item(5, 11).Value = 1.056900152086484
$ws.Cells.Item(5, 12).Value = 0.3513574990490724

$ws.Cells.Item(6, 2).Value = 1.202997182580049
$ws.Cells.Item(6, 3).Value = 0.06245003754649758
$ws.Cells.Item(6, 4).Value = 0.004277372990814765
$ws.Cells.Item(6, 5).Value = 0.06562834507854376
$ws.Cells.Item(6, 6).Value = 4.565733381866977
$ws.Cells.Item(6, 8).Value = 0.07973214163530429
$ws.Cells.Item(6, 9).Value = 2.962330757970662
$ws.Cells.Item(6, 10).Value = 0.1737082563564059
$ws.Cells.Item(6, 11).Value = 1.056050769841846
$ws.Cells.Item(6, 12).Value = 0.3513124926112354

$ws.Cells.Item(7, 2).Value = 1.207202041390588
$ws.Cells.Item(7, 3).Value = 0.06393852312737636
$ws.Cells.Item(7, 4).Value = 0.004233872854239706
$ws.Cells.Item(7, 5).Value = 0.06554188825680729
$ws.Cells.Item(7, 6).Value = 4.576283847233796
$ws.Cells.Item(7, 8).Value = 0.07973214163530429
$ws.Cells.Item(7, 9).Value = 2.968839543011413
$ws.Cells.Item(7, 10).Value = 0.1737112624479273
$ws.Cells.Item(7, 11).Value = 1.062015252128248
$ws.Cells.Item(7, 12).Value = 0.3516406043675744

$ws.Cells.Item(8, 2).Value = 1.227619916911664
$ws.Cells.Item(8, 3).Value = 0.07060472970725584
$ws.Cells.Item(8, 4).Value = 0.004057136277395124
$ws.Cells.Item(8, 5).Value = 0.06520664220990113
$ws.Cells.Item(8, 6).Value = 4.625619690632021
$ws.Cells.Item(8, 8).Value = 0.07973214163530429
$ws.Cells.Item(8, 9).Value = 2.999082410949214
$ws.Cells.Item(8, 10).Value = 0.1737980450870289
$ws.Cells.Item(8, 11).Value = 1.090038064358765
$ws.Cells.Item(8, 12).Value = 0.3534770064474202

$ws.Cells.Item(9, 2).Value = 1.274341815827171
$ws.Cells.Item(9, 3).Value = 0.08403464828941765
$ws.Cells.Item(9, 4).Value = 0.003764836723097353
$ws.Cells.Item(9, 5).Value = 0.06471458456161727
$ws.Cells.Item(9, 6).Value = 4.732381880448486
$ws.Cells.Item(9, 8).Value = 0.07973214163530429
$ws.Cells.Item(9, 9).Value = 3.063876326221433
$ws.Cells.Item(9, 10).Value = 0.1742320158728248
$ws.Cells.Item(9, 11).Value = 1.151110071799508
$ws.Cells.Item(9, 12).Value = 0.3584712324411328

$ws.Cells.Item(10, 2).Value = 1.312739641244917
$ws.Cells.Item(10, 3).Value = 0.09412641534504473
$ws.Cells.Item(10, 4).Value = 0.003582800054517943
$ws.Cells.Item(10, 5).Value = 0.06445307617095075
$ws.Cells.Item(10, 6).Value = 4.816965553773628
$ws.Cells.Item(10, 8).Value = 0.07973214163530429
$ws.Cells.Item(10, 9).Value = 3.114854153546688
$ws.Cells.Item(10, 10).Value = 0.1747114558830383
$ws.Cells.Item(10, 11).Value = 1.199719168271145
$ws.Cells.Item(10, 12).Value = 0.3629880283305056

$ws.Cells.Item(11, 2).Value = 1.331091868551681
$ws.Cells.Item(11, 3).Value = 0.09876778988007118
$ws.Cells.Item(11, 4).Value = 0.003507012920325625
$ws.Cells.Item(11, 5).Value = 0.06435565919585784
$ws.Cells.Item(11, 6).Value = 4.856788753741057
$ws.Cells.Item(11, 8).Value = 0.07973214163530429
$ws.Cells.Item(11, 9).Value = 3.138784384221495
$ws.Cells.Item(11, 10).Value = 0.1749644055662323
$ws.Cells.Item(11, 11).Value = 1.222647243526183
$ws.Cells.Item(11, 12).Value = 0.3652265989876042

$ws.Cells.Item(12, 2).Value = 1.338168507203477
$ws.Cells.Item(12, 3).Value = 0.1005327346339016
$ws.Cells.Item(12, 4).Value = 0.003479317904596879
$ws.Cells.Item(12, 5).Value = 0.06432185417751235
$ws.Cells.Item(12, 6).Value = 4.872062799142839
$ws.Cells.Item(12, 8).Value = 0.07973214163530429
$ws.Cells.Item(12, 9).Value = 3.147952963583094
$ws.Cells.Item(12, 10).Value = 0.1750651971897668
$ws.Cells.Item(12, 11).Value = 1.231446856668356
$ws.Cells.Item(12, 12).Value = 0.3661006867057637

$ws.Cells.Item(13, 2).Value = 1.336638778933178
$ws.Cells.Item(13, 3).Value = 0.1001522936565777
$ws.Cells.Item(13, 4).Value = 0.003485237965063614
$ws.Cells.Item(13, 5).Value = 0.06432899771318823
$ws.Cells.Item(13, 6).Value = 4.868764632249736
$ws.Cells.Item(13, 8).Value = 0.07973214163530429
$ws.Cells.Item(13, 9).Value = 3.145973594973043
$ws.Cells.Item(13, 10).Value = 0.1750432674079221
$ws.Cells.Item(13, 11).Value = 1.229546488205045
$ws.Cells.Item(13, 12).Value = 0.3659112633563524

$ws.Cells.Item(14, 2).Value = 1.331671522712639
$ws.Cells.Item(14, 3).Value = 0.0989128450200667
$ws.Cells.Item(14, 4).Value = 0.0035047143486171
$ws.Cells.Item(14, 5).Value = 0.06435281629134693
$ws.Cells.Item(14, 6).Value = 4.85804147125387
$ws.Cells.Item(14, 8).Value = 0.07973214163530429
$ws.Cells.Item(14, 9).Value = 3.139536548283843
$ws.Cells.Item(14, 10).Value = 0.1749725974802701
$ws.Cells.Item(14, 11).Value = 1.223368843281719
$ws.Cells.Item(14, 12).Value = 0.3652979820990652

$ws.Cells.Item(15, 2).Value = 1.328645474703109
$ws.Cells.Item(15, 3).Value = 0.09815460753551974
$ws.Cells.Item(15, 4).Value = 0.003516774763741015
$ws.Cells.Item(15, 5).Value = 0.06436780717882584
$ws.Cells.Item(15, 6).Value = 4.851498483060141
$ws.Cells.Item(15, 8).Value = 0.07973214163530429
$ws.Cells.Item(15, 9).Value = 3.135607577983635
$ws.Cells.Item(15, 10).Value = 0.1749299617313973
$ws.Cells.Item(15, 11).Value = 1.21960012597205
$ws.Cells.Item(15, 12).Value = 0.3649257648963982

$ws.Cells.Item(16, 2).Value = 1.311558068985732
$ws.Cells.Item(16, 3).Value = 0.09382411551999326
$ws.Cells.Item(16, 4).Value = 0.003587893723773661
$ws.Cells.Item(16, 5).Value = 0.06445987498169004
$ws.Cells.Item(16, 6).Value = 4.814390122761097
$ws.Cells.Item(16, 8).Value = 0.07973214163530429
$ws.Cells.Item(16, 9).Value = 3.11330517503913
$ws.Cells.Item(16, 10).Value = 0.1746956257544419
$ws.Cells.Item(16, 11).Value = 1.198237173153245
$ws.Cells.Item(16, 12).Value = 0.3628454279469224

$ws.Cells.Item(17, 2).Value = 1.301301993251485
$ws.Cells.Item(17, 3).Value = 0.09118051272575656
$ws.Cells.Item(17, 4).Value = 0.003633317271008796
$ws.Cells.Item(17, 5).Value = 0.06452186408597171
$ws.Cells.Item(17, 6).Value = 4.79197025609588
$ws.Cells.Item(17, 8).Value = 0.07973214163530429
$ws.Cells.Item(17, 9).Value = 3.099813149183973
$ws.Cells.Item(17, 10).Value = 0.1745607896129435
$ws.Cells.Item(17, 11).Value = 1.18534055446483
$ws.Cells.Item(17, 12).Value = 0.3616162637167406

$ws.Cells.Item(18, 2).Value = 1.295486268410741
$ws.Cells.Item(18, 3).Value = 0.08966473799732455
$ws.Cells.Item(18, 4).Value = 0.003660105008146175
$ws.Cells.Item(18, 5).Value = 0.06455954724298074
$ws.Cells.Item(18, 6).Value = 4.779201580168262
$ws.Cells.Item(18, 8).Value = 0.07973214163530429
$ws.Cells.Item(18, 9).Value = 3.092122543051715
$ws.Cells.Item(18, 10).Value = 0.1744865160507842
$ws.Cells.Item(18, 11).Value = 1.177999525345001
$ws.Cells.Item(18, 12).Value = 0.3609265893473861

$ws.Cells.Item(19, 2).Value = 1.293531476572639
$ws.Cells.Item(19, 3).Value = 0.08915233697466363
$ws.Cells.Item(19, 4).Value = 0.003669288635847145
$ws.Cells.Item(19, 5).Value = 0.06457265494765174
$ws.Cells.Item(19, 6).Value = 4.774900065283759
$ws.Cells.Item(19, 8).Value = 0.07973214163530429
$ws.Cells.Item(19, 9).Value = 3.089530589768017
$ws.Cells.Item(19, 10).Value = 0.1744619319426839
$ws.Cells.Item(19, 11).Value = 1.17552716441287
$ws.Cells.Item(19, 12).Value = 0.3606960516987812

$ws.Cells.Item(20, 2).Value = 1.302385150372317
$ws.Cells.Item(20, 3).Value = 0.09146143595664569
$ws.Cells.Item(20, 4).Value = 0.003628413456659541
$ws.Cells.Item(20, 5).Value = 0.06451505537849833
$ws.Cells.Item(20, 6).Value = 4.794343780042908
$ws.Cells.Item(20, 8).Value = 0.07973214163530429
$ws.Cells.Item(20, 9).Value = 3.101242187791655
$ws.Cells.Item(20, 10).Value = 0.1745748036652373
$ws.Cells.Item(20, 11).Value = 1.186705476903285
$ws.Cells.Item(20, 12).Value = 0.3617453193548243

$ws.Cells.Item(21, 2).Value = 1.333127079716604
$ws.Cells.Item(21, 3).Value = 0.09927670079184736
$ws.Cells.Item(21, 4).Value = 0.003498966461842112
$ws.Cells.Item(21, 5).Value = 0.06434573658655118
$ws.Cells.Item(21, 6).Value = 4.861185859273007
$ws.Cells.Item(21, 8).Value = 0.07973214163530429
$ws.Cells.Item(21, 9).Value = 3.141424366282365
$ws.Cells.Item(21, 10).Value = 0.1749932191649464
$ws.Cells.Item(21, 11).Value = 1.225180185422062
$ws.Cells.Item(21, 12).Value = 0.365477401847258

$ws.Cells.Item(22, 2).Value = 1.353959159478222
$ws.Cells.Item(22, 3).Value = 0.1044273496265191
$ws.Cells.Item(22, 4).Value = 0.003420214366437246
$ws.Cells.Item(22, 5).Value = 0.06425305069642029
$ws.Cells.Item(22, 6).Value = 4.906001179382969
$ws.Cells.Item(22, 8).Value = 0.07973214163530429
$ws.Cells.Item(22, 9).Value = 3.168307924928357
$ws.Cells.Item(22, 10).Value = 0.17529584840225
$ws.Cells.Item(22, 11).Value = 1.251009040047677
$ws.Cells.Item(22, 12).Value = 0.3680703305260238

$ws.Cells.Item(23, 2).Value = 1.342772996903079
$ws.Cells.Item(23, 3).Value = 0.1016743961668283
$ws.Cells.Item(23, 4).Value = 0.003461712572167208
$ws.Cells.Item(23, 5).Value = 0.06430087869505563
$ws.Cells.Item(23, 6).Value = 4.881978870139619
$ws.Cells.Item(23, 8).Value = 0.07973214163530429
$ws.Cells.Item(23, 9).Value = 3.153902638996385
$ws.Cells.Item(23, 10).Value = 0.1751316624355184
$ws.Cells.Item(23, 11).Value = 1.237161168987313
$ws.Cells.Item(23, 12).Value = 0.3666723790548758

$ws.Cells.Item(24, 2).Value = 1.301895203872704
$ws.Cells.Item(24, 3).Value = 0.09133441790075381
$ws.Cells.Item(24, 4).Value = 0.003630628374399869
$ws.Cells.Item(24, 5).Value = 0.06451812722511363
$ws.Cells.Item(24, 6).Value = 4.793270333592432
$ws.Cells.Item(24, 8).Value = 0.07973214163530429
$ws.Cells.Item(24, 9).Value = 3.10059591343466
$ws.Cells.Item(24, 10).Value = 0.1745684578021987
$ws.Cells.Item(24, 11).Value = 1.186088166809981
$ws.Cells.Item(24, 12).Value = 0.3616869203812882

$ws.Cells.Item(25, 2).Value = 1.260987195016099
$ws.Cells.Item(25, 3).Value = 0.08036252016950129
$ws.Cells.Item(25, 4).Value = 0.003838139503849369
$ws.Cells.Item(25, 5).Value = 0.06483007966772902
$ws.Cells.Item(25, 6).Value = 4.702423590695986
$ws.Cells.Item(25, 8).Value = 0.07973214163530429
$ws.Cells.Item(25, 9).Value = 3.045757439977592
$ws.Cells.Item(25, 10).Value = 0.1740863970475672
$ws.Cells.Item(25, 11).Value = 1.133932460600192
$ws.Cells.Item(25, 12).Value = 0.3569712075155991
